$wb = $excel.ActiveWorkbook

# Rename sheets to match 'pregnant women'
$wsAffFrac = $wb.Worksheets.Item("Interventions maternal aff frac")
$wsAffFrac.Name = "Inter. pregnant women aff frac"

$wsEff = $wb.Worksheets.Item("Interventions maternal eff")
$wsEff.Name = "Inter. pregnant women eff"

# Alter test mortality rates for pregnant women on the "mortality rates" sheet
$wsMort = $wb.Worksheets.Item("mortality rates")
$wsMort.Range("E2:J2").Value = 0.01

# Update the active sheet / selection state
$wsMort.Activate()
$wsMort.Range("E2:J2").Select()
